# Added command to drop rows that are NA.
# To exercise / demonstrate the new "drop NA rows" behaviour, two extra
# sample rows are appended to the bottom of the data:
#   - row 5: an NNSS value that is just a blank space (" ") -> treated as NA
#   - row 6: a normal, valid NNSS value (7560000000004)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - a "blank"/NA-like entry in column A (NNSS)
$ws.Cells.Item(5, 1).Value = " "

# Row 6 - a normal NNSS id, formatted the same way as the existing ids (A2:A4)
$ws.Cells.Item(6, 1).Value = 7560000000004
$ws.Cells.Item(6, 1).Style = $ws.Cells.Item(4, 1).Style
$ws.Range("A6").NumberFormat = $ws.Range("A4").NumberFormat

# Move the active selection, matching where the author's cursor ended up
$ws.Range("D14").Select()
